$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($r1, $r2) {
    $rng1 = $ws.Range("B$r1`:AD$r1")
    $rng2 = $ws.Range("B$r2`:AD$r2")
    $v1 = $rng1.Value()
    $v2 = $rng2.Value()
    $rng1.Value = $v2
    $rng2.Value = $v1
}

Swap-Rows 20 21
Swap-Rows 59 60
Swap-Rows 215 216
Swap-Rows 226 227
Swap-Rows 252 253
Swap-Rows 271 272
